$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.184.02"
$ws.Range("D3").Value = "1.911.08"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("D5").Value = "'314.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("E7").Value = "  +0.57%  "
$ws.Range("D8").Value = "'0.3923"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").Value = "'0.09300"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.54%  "
$ws.Range("D10").Value = "'1.140"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.34%  "
$ws.Range("D11").Value = "'41.89"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.34%  "
$ws.Range("D12").Value = "'6.394"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.58%  "
$ws.Range("D13").Value = "'20.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.75%  "
$ws.Range("D14").Value = "1.907.21"
$ws.Range("E14").Value = "  +1.43%  "
$ws.Range("D15").Value = "'7.314"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.61%  "
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").Value = "'0.00001120"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("D18").Value = "'92.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("D19").Value = "'0.06608"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "'17.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.89%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("E22").Value = "  +0.63%  "
$ws.Range("D23").Value = "28.229.26"
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("D24").Value = "'11.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.65%  "
$ws.Range("E25").Value = "  +1.42%  "
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").Value = "'3.397"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "'2.590"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.78%  "
$ws.Range("B28").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C28").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D28").Value = "2.122.56"
$ws.Range("E28").Value = "  +1.30%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'21.10"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.85%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "'158.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.37%  "
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Value = "'127.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'1.096"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.45%  "
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").Value = "'0.1076"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.69%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'5.636"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'3.616"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = "'9.702"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.02%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.06667"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.04%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02429"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.05%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "'1.244"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "'0.2198"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'1.289"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.98%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.6458"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.33%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "'11.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("B44").Value = "InternetComputer(DFINITY)"
$ws.Range("C44").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D44").Value = "'4.997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "'1.000"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'13.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.44%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.6053"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.93%  "
$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").Value = "'3.721"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("B49").Value = "WEMIXTOKEN"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "'1.283"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.31%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'2.013"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "'123.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.58%  "

Write-Output "done"
